# Scheduled runner: refresh Leve profit-calc market data (currentAveragePrice*, LevePrice*, LeveProfit*)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 329
$ws.Range("I2").Value = 96
$ws.Range("J2").Value = 678.5
$ws.Range("K2").Value = 96
$ws.Range("L2").Value = 678.5
$ws.Range("M2").Value = 17
$ws.Range("N2").Value = -904.5
$ws.Range("H40").Value = 3999.6191
$ws.Range("I40").Value = 3997.3333
$ws.Range("K40").Value = 3997.3333
$ws.Range("M40").Value = -3822.3333
$ws.Range("H62").Value = 7933.647
$ws.Range("I62").Value = 7345.533
$ws.Range("K62").Value = 7345.533
$ws.Range("M62").Value = -6721.533
$ws.Range("H65").Value = 7933.647
$ws.Range("I65").Value = 7345.533
$ws.Range("K65").Value = 36727.665
$ws.Range("M65").Value = -33607.665
$ws.Range("H88").Value = 1942.8
$ws.Range("J88").Value = 1942.8
$ws.Range("L88").Value = 1942.8
$ws.Range("N88").Value = -2754.8
$ws.Range("H91").Value = 1942.8
$ws.Range("J91").Value = 1942.8
$ws.Range("L91").Value = 1942.8
$ws.Range("N91").Value = -4750.8
$ws.Range("H98").Value = 1383.9642
$ws.Range("I98").Value = 1328.579
$ws.Range("J98").Value = 1500.8889
$ws.Range("K98").Value = 1328.579
$ws.Range("L98").Value = 1500.8889
$ws.Range("M98").Value = 169.421
$ws.Range("N98").Value = -4496.8889
$ws.Range("H100").Value = 2462.4119
$ws.Range("H107").Value = 2117.2856
$ws.Range("J107").Value = 432.66666
$ws.Range("L107").Value = 432.66666
$ws.Range("N107").Value = -4272.66666
$ws.Range("H113").Value = 4996.575
$ws.Range("I113").Value = 4842.727
$ws.Range("K113").Value = 4842.727
$ws.Range("M113").Value = -1588.727
$ws.Range("H122").Value = 1383.9642
$ws.Range("I122").Value = 1328.579
$ws.Range("J122").Value = 1500.8889
$ws.Range("K122").Value = 3985.737
$ws.Range("L122").Value = 4502.6667
$ws.Range("M122").Value = -1535.737
$ws.Range("N122").Value = -9402.6667
$ws.Range("H135").Value = 1568.8667
$ws.Range("I135").Value = 1653.7693
$ws.Range("J135").Value = 1017
$ws.Range("K135").Value = 14883.9237
$ws.Range("L135").Value = 9153
$ws.Range("M135").Value = -12348.9237
$ws.Range("N135").Value = -14223

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2192.077
$ws.Range("I2").Value = 1314
$ws.Range("J2").Value = 2740.875
$ws.Range("K2").Value = 1314
$ws.Range("L2").Value = 2740.875
$ws.Range("M2").Value = -1201
$ws.Range("N2").Value = -2966.875
$ws.Range("H45").Value = 3810.2163
$ws.Range("I45").Value = 3254.7058
$ws.Range("K45").Value = 3254.7058
$ws.Range("M45").Value = -2877.7058
$ws.Range("H97").Value = 1455.8
$ws.Range("I97").Value = 1187
$ws.Range("J97").Value = 2195
$ws.Range("K97").Value = 1187
$ws.Range("L97").Value = 2195
$ws.Range("M97").Value = -691
$ws.Range("N97").Value = -3187
$ws.Range("H107").Value = 34994
$ws.Range("J107").Value = 34994
$ws.Range("L107").Value = 34994
$ws.Range("N107").Value = -42674
$ws.Range("H108").Value = 68659.336
$ws.Range("J108").Value = 68659.336
$ws.Range("L108").Value = 68659.336
$ws.Range("N108").Value = -76339.336
$ws.Range("H109").Value = 55613.332
$ws.Range("J109").Value = 55613.332
$ws.Range("L109").Value = 55613.332
$ws.Range("N109").Value = -58387.332
$ws.Range("H110").Value = 1913.75
$ws.Range("I110").Value = 1301.6666
$ws.Range("K110").Value = 1301.6666
$ws.Range("M110").Value = 743.3334
$ws.Range("H113").Value = 95000
$ws.Range("J113").Value = 95000
$ws.Range("L113").Value = 95000
$ws.Range("N113").Value = -103678
$ws.Range("H115").Value = 94493.5
$ws.Range("J115").Value = 94493.5
$ws.Range("L115").Value = 94493.5
$ws.Range("N115").Value = -97627.5
$ws.Range("H116").Value = 2192.077
$ws.Range("I116").Value = 1314
$ws.Range("J116").Value = 2740.875
$ws.Range("K116").Value = 1314
$ws.Range("L116").Value = 2740.875
$ws.Range("M116").Value = 980
$ws.Range("N116").Value = -7328.875
$ws.Range("H118").Value = 126988
$ws.Range("J118").Value = 126988
$ws.Range("L118").Value = 126988
$ws.Range("N118").Value = -130302
$ws.Range("H119").Value = 80601.336
$ws.Range("J119").Value = 80601.336
$ws.Range("L119").Value = 80601.336
$ws.Range("N119").Value = -90277.336
$ws.Range("H120").Value = 67750.836
$ws.Range("J120").Value = 67750.836
$ws.Range("L120").Value = 67750.836
$ws.Range("N120").Value = -77426.836
$ws.Range("H121").Value = 114819
$ws.Range("J121").Value = 114819
$ws.Range("L121").Value = 114819
$ws.Range("N121").Value = -118313
$ws.Range("H122").Value = 4441.9
$ws.Range("I122").Value = 3837.9167
$ws.Range("K122").Value = 11513.7501
$ws.Range("M122").Value = -9063.750100000001
$ws.Range("H132").Value = 3205.4255
$ws.Range("I132").Value = 3006.8108
$ws.Range("J132").Value = 3940.3
$ws.Range("K132").Value = 9020.432400000002
$ws.Range("L132").Value = 11820.9
$ws.Range("M132").Value = -6490.432400000002
$ws.Range("N132").Value = -16880.9

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2192.077
$ws.Range("I3").Value = 1314
$ws.Range("J3").Value = 2740.875
$ws.Range("K3").Value = 1314
$ws.Range("L3").Value = 2740.875
$ws.Range("M3").Value = -1200
$ws.Range("N3").Value = -2968.875
$ws.Range("H22").Value = 476.57144
$ws.Range("I22").Value = 272
$ws.Range("K22").Value = 272
$ws.Range("M22").Value = -99
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7501.5293
$ws.Range("I31").Value = 3346.8333
$ws.Range("K31").Value = 3346.8333
$ws.Range("M31").Value = -3051.8333
$ws.Range("H34").Value = 7501.5293
$ws.Range("I34").Value = 3346.8333
$ws.Range("K34").Value = 3346.8333
$ws.Range("M34").Value = -3144.8333
$ws.Range("H99").Value = 3692
$ws.Range("I99").Value = 3692
$ws.Range("K99").Value = 3692
$ws.Range("M99").Value = -2194
$ws.Range("H126").Value = 3692
$ws.Range("I126").Value = 3692
$ws.Range("K126").Value = 11076
$ws.Range("M126").Value = -8606
$ws.Range("H132").Value = 3550.7273
$ws.Range("I132").Value = 3338.2222
$ws.Range("J132").Value = 4507
$ws.Range("K132").Value = 10014.6666
$ws.Range("L132").Value = 13521
$ws.Range("M132").Value = -7484.6666
$ws.Range("N132").Value = -18581
$ws.Range("H134").Value = 1979
$ws.Range("I134").Value = 1977.6923
$ws.Range("K134").Value = 5933.0769
$ws.Range("M134").Value = -3398.0769

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 526377.1
$ws.Range("J12").Value = 625068.1
$ws.Range("L12").Value = 1875204.3
$ws.Range("N12").Value = -1875550.3
$ws.Range("H107").Value = 643.1852
$ws.Range("I107").Value = 426.29413
$ws.Range("K107").Value = 1278.88239
$ws.Range("M107").Value = 641.11761

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1075.25
$ws.Range("I97").Value = 562.44446
$ws.Range("K97").Value = 562.44446
$ws.Range("M97").Value = -66.44446000000005
$ws.Range("H103").Value = 95346
$ws.Range("J103").Value = 95346
$ws.Range("L103").Value = 95346
$ws.Range("N103").Value = -97690
$ws.Range("H132").Value = 2701.6843
$ws.Range("I132").Value = 2576.1428
$ws.Range("J132").Value = 4166.3335
$ws.Range("K132").Value = 7728.428400000001
$ws.Range("L132").Value = 12499.0005
$ws.Range("M132").Value = -5198.428400000001
$ws.Range("N132").Value = -17559.0005

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3667.0454
$ws.Range("J46").Value = 3760.7144
$ws.Range("L46").Value = 3760.7144
$ws.Range("N46").Value = -4136.7144
$ws.Range("H119").Value = 99925
$ws.Range("J119").Value = 99925
$ws.Range("L119").Value = 99925
$ws.Range("N119").Value = -109601

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 60000
$ws.Range("J99").Value = 60000
$ws.Range("L99").Value = 60000
$ws.Range("N99").Value = -65990

